$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing row 41, column C: 清涼地 -> 清涼堂
$ws.Range("C41").Value = "第91期 秘寶 開放區域 清涼堂 祕寶效果: 透過元素袋獲取星途解讀道具有5%(18%) 翻倍"

# New row 42
$ws.Range("A42").Value = "2026/1/9"
$ws.Range("B42").Value = "2026/3/6"
$ws.Range("C42").Value = "第92期 第四代寵物"

# New row 43
$ws.Range("A43").Value = "2026/1/16"
$ws.Range("B43").Value = "2026/3/13"
$ws.Range("C43").Value = "第93期 秘寶 開放區域 極樂島 祕寶效果: ??"

# New row 44
$ws.Range("A44").Value = "2026/1/23"
$ws.Range("B44").Value = "2026/3/20"
$ws.Range("C44").Value = "第94期 第五代寵物"

# Match the new selection/active cell shown in the diff
$ws.Range("C50").Select()
